$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update F column (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 315
$ws1.Range("F4").Value = 2838
$ws1.Range("F6").Value = 601

# Sheet "全部类型" (all types) - update F column (想去人数 / interested count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 315
$ws4.Range("F6").Value = 2838
$ws4.Range("F8").Value = 601
